$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (LinearRegression) - only C2/D2 change
$ws.Range("C2").Value = 0.08176149744867572
$ws.Range("D2").Value = 0.08176149744867572

# Row 3 (RandomForestRegressor)
$ws.Range("B3").Value = 0.02221210674801416
$ws.Range("C3").Value = 0.02232892104023649
$ws.Range("D3").Value = 0.02222514405453239

# Row 4: label change GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.02341432668715818
$ws.Range("C4").Value = 0.02213485332120204
$ws.Range("D4").Value = 0.02195063456506316

# Row 5: label change AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.01848119518831166
$ws.Range("C5").Value = 0.01828038281338033
$ws.Range("D5").Value = 0.01916125283641585
